$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 ("H 72" record); this shifts all subsequent rows up by one
# and updates the used range from A1:F63 to A1:F62.
$ws.Rows.Item(2).Delete()
